# Refresh the cryptos worksheet with the latest coinranking.com snapshot.
# Each hashtable below represents one spreadsheet row that changed since the
# last run: rank-shuffles touch Coin/Link/Price, pure price ticks only touch
# Price/Volume(1h).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="26.597.53"; E="  +4.08%  " }
    @{ Row=3; D="1.743.02"; E="  +4.42%  " }
    @{ Row=4; E="  +0.04%  " }
    @{ Row=5; D="246.43"; E="  +3.90%  " }
    @{ Row=6; D="1.000"; E="  -0.01%  " }
    @{ Row=7; D="0.4825"; E="  +1.81%  " }
    @{ Row=8; D="0.2692"; E="  +3.49%  " }
    @{ Row=9; D="0.06264"; E="  +1.55%  " }
    @{ Row=10; D="1.745.65"; E="  +4.59%  " }
    @{ Row=11; D="0.07144"; E="  +2.02%  " }
    @{ Row=12; D="15.92"; E="  +7.70%  " }
    @{ Row=13; D="0.6250"; E="  +6.60%  " }
    @{ Row=14; D="4.522"; E="  +3.61%  " }
    @{ Row=15; D="77.47"; E="  +2.70%  " }
    @{ Row=16; E="  -0.03%  " }
    @{ Row=17; D="26.602.55"; E="  +4.13%  " }
    @{ Row=18; D="1.000"; E="  +0.09%  " }
    @{ Row=19; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="11.81"; E="  +3.53%  " }
    @{ Row=20; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000006908"; E="  +2.69%  " }
    @{ Row=21; D="1.967.90"; E="  +4.53%  " }
    @{ Row=22; D="4.622"; E="  +4.07%  " }
    @{ Row=23; D="8.886"; E="  +1.29%  " }
    @{ Row=24; D="5.368"; E="  +2.69%  " }
    @{ Row=25; D="136.28"; E="  -0.60%  " }
    @{ Row=27; D="1.815"; E="  +5.62%  " }
    @{ Row=28; D="1.430"; E="  +2.99%  " }
    @{ Row=29; D="106.81"; E="  +2.31%  " }
    @{ Row=30; D="4.016"; E="  +0.34%  " }
    @{ Row=31; D="3.746"; E="  +3.33%  " }
    @{ Row=32; D="0.07888"; E="  +0.74%  " }
    @{ Row=33; D="0.04594"; E="  +6.62%  " }
    @{ Row=34; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="0.9993"; E="  +0.02%  " }
    @{ Row=35; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.617"; E="  -0.29%  " }
    @{ Row=36; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.6415"; E="  +5.92%  " }
    @{ Row=37; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.000"; E="  +4.77%  " }
    @{ Row=38; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.9314"; E="  -0.29%  " }
    @{ Row=39; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="113.95"; E="  +14.61%  " }
    @{ Row=40; D="2.437"; E="  -3.42%  " }
    @{ Row=41; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.987"; E="  +7.23%  " }
    @{ Row=42; D="5.797"; E="  +18.34%  " }
    @{ Row=43; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.000"; E="  +0.03%  " }
    @{ Row=44; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01512"; E="  +2.25%  " }
    @{ Row=45; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.3924"; E="  +4.67%  " }
    @{ Row=46; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1218"; E="  +9.15%  " }
    @{ Row=47; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="6.755"; E="  +8.81%  " }
    @{ Row=48; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.05337"; E="  +1.44%  " }
    @{ Row=49; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="7.933"; E="  +6.00%  " }
    @{ Row=50; B="Elrond"; C="https://coinranking.com/coin/omwkOTglq+elrond-egld"; D="30.81"; E="  +3.13%  " }
    @{ Row=51; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="1.266"; E="  +5.22%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Price column holds free-form text (e.g. "26.597.53",
        # "1.000"), not a number -- force text so Excel does not
        # reinterpret/round it, then drop the format change so the
        # cell style matches the original (unstyled) cell.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $u.E }
}

